$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Purchase Price (INR) and Current Price (INR) totals to reflect profit booking
$ws.Range("B2").Value = 233673.2
$ws.Range("C2").Value = 258116.7500762939

# Recalculate derived Gain/Loss (INR) and Gain/Loss (%) totals
$b2 = $ws.Range("B2").Value()
$c2 = $ws.Range("C2").Value()
$d2 = $c2 - $b2
$e2 = ($d2 / $b2) * 100

$ws.Range("D2").Value = $d2
$ws.Range("E2").Value = $e2
